# Updates cryptos list price (D) and volume/1h (E) columns, and fixes the
# row order/content for Kaspa/Stacks (rows 38-39), per the scraped commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Price($row, $price) {
    # Many of the price strings (e.g. "594.89") look like genuine numbers to
    # Excel and would otherwise be auto-converted to a numeric cell (losing
    # the exact textual representation, e.g. trailing zeros). Force the cell
    # to Text format while assigning, then restore the default ("Normal")
    # style so we don't leave stray per-cell number formatting behind.
    $c = $ws.Range("D$row")
    $c.NumberFormat = "@"
    $c.Value = $price
    $c.Style = "Normal"
}

function Set-Volume($row, $volume) {
    $ws.Range("E$row").Value = $volume
}

Set-Price  2  "61.387.55"
Set-Volume 2  "  +0.60%  "

Set-Price  3  "2.929.37"
Set-Volume 3  "  +0.23%  "

Set-Volume 4  "  +0.06%  "

Set-Price  5  "594.89"
Set-Volume 5  "  +0.80%  "

Set-Price  6  "144.90"
Set-Volume 6  "  -0.21%  "

Set-Volume 7  "  +0.08%  "

Set-Volume 8  "  -0.75%  "

Set-Price  9  "6.97"
Set-Volume 9  "  +1.61%  "

Set-Volume 10 "  -1.61%  "

Set-Price  11 "0.439"
Set-Volume 11 "  -0.70%  "

Set-Volume 12 "  -0.76%  "

Set-Price  13 "33.62"
Set-Volume 13 "  -0.07%  "

Set-Volume 14 "  +0.58%  "

Set-Price  15 "3.418.46"
Set-Volume 15 "  +0.39%  "

Set-Price  16 "61.438.36"
Set-Volume 16 "  +0.70%  "

Set-Price  17 "6.71"
Set-Volume 17 "  -0.14%  "

Set-Price  18 "2.932.11"
Set-Volume 18 "  +0.41%  "

Set-Price  19 "432.07"
Set-Volume 19 "  +0.28%  "

Set-Price  20 "13.48"
Set-Volume 20 "  -0.10%  "

Set-Price  21 "0.678"
Set-Volume 21 "  -0.88%  "

Set-Price  22 "7.11"
Set-Volume 22 "  +0.24%  "

Set-Price  23 "81.88"
Set-Volume 23 "  +0.97%  "

Set-Price  24 "10.90"
Set-Volume 24 "  -1.50%  "

Set-Price  25 "2.19"
Set-Volume 25 "  -1.30%  "

Set-Price  26 "11.77"
Set-Volume 26 "  -2.32%  "

Set-Volume 27 "  -0.06%  "

Set-Volume 28 "  -3.49%  "

Set-Price  29 "2.61"
Set-Volume 29 "  -0.44%  "

Set-Price  30 "6.92"
Set-Volume 30 "  -2.73%  "

Set-Volume 31 "  +0.70%  "

Set-Volume 32 "  +1.14%  "

Set-Volume 33 "  +0.02%  "

Set-Price  34 "0.0₃0878"
Set-Volume 34 "  +1.95%  "

Set-Volume 35 "  +0.12%  "

Set-Price  36 "5.64"
Set-Volume 36 "  +0.02%  "

Set-Price  37 "3.00"
Set-Volume 37 "  -1.97%  "

# Rows 38 and 39 swap identity: Kaspa/Stacks order flips with refreshed data.
$ws.Range("B38").Value = "Stacks"
$ws.Range("C38").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-Price  38 "2.00"
Set-Volume 38 "  +0.03%  "

$ws.Range("B39").Value = "Kaspa"
$ws.Range("C39").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-Price  39 "0.124"
Set-Volume 39 "  -0.94%  "

Set-Volume 40 "  -0.02%  "

Set-Price  41 "42.54"
Set-Volume 41 "  +7.96%  "

Set-Price  42 "0.282"
Set-Volume 42 "  -1.35%  "

Set-Price  43 "0.0347"
Set-Volume 43 "  -0.10%  "

Set-Price  44 "2.703.63"
Set-Volume 44 "  -0.21%  "

Set-Price  45 "134.41"
Set-Volume 45 "  +2.10%  "

Set-Price  46 "364.69"
Set-Volume 46 "  -3.01%  "

Set-Price  48 "23.71"
Set-Volume 48 "  -1.90%  "

Set-Volume 49 "  -1.36%  "

Set-Price  50 "2.00"
Set-Volume 50 "  -1.93%  "

Set-Volume 51 "  -0.56%  "
